# excess_mortality_provinces.xlsx - "[2022-07-26] - Monkeypox update"
# Apply the weekly-data refresh: a handful of previously published weeks
# were corrected (minor upward revisions to several province counts),
# and two new ISO weeks (2022 week 27 and 2022 week 28) were appended
# with their source counts and the derived excess-mortality percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Revised weekly figures (minor data corrections) ---
$ws.Range("Y111").Value = 80
$ws.Range("AA112").Value = 241
$ws.Range("X113").Value = 689
$ws.Range("W115").Value = 500
$ws.Range("X115").Value = 623
$ws.Range("X116").Value = 678
$ws.Range("W117").Value = 516
$ws.Range("W118").Value = 530
$ws.Range("W120").Value = 545
$ws.Range("X120").Value = 731
$ws.Range("V121").Value = 236
$ws.Range("W121").Value = 558
$ws.Range("Z121").Value = 503
$ws.Range("AA121").Value = 304
$ws.Range("X124").Value = 637
$ws.Range("V125").Value = 208
$ws.Range("W125").Value = 514
$ws.Range("AA125").Value = 228
$ws.Range("X126").Value = 571
$ws.Range("X127").Value = 589
$ws.Range("P128").Value = 107
$ws.Range("T128").Value = 53
$ws.Range("U128").Value = 353
$ws.Range("W128").Value = 455
$ws.Range("X128").Value = 624
$ws.Range("U129").Value = 395
$ws.Range("W129").Value = 438
$ws.Range("X129").Value = 589
$ws.Range("P130").Value = 115
$ws.Range("T130").Value = 44
$ws.Range("V130").Value = 221
$ws.Range("W130").Value = 446
$ws.Range("X130").Value = 605
$ws.Range("Z130").Value = 466
$ws.Range("AA130").Value = 227
$ws.Range("P131").Value = 116
$ws.Range("U131").Value = 383
$ws.Range("V131").Value = 210
$ws.Range("W131").Value = 440
$ws.Range("X131").Value = 697
$ws.Range("Z131").Value = 451
$ws.Range("AA131").Value = 211
$ws.Range("P132").Value = 116
$ws.Range("Q132").Value = 149
$ws.Range("S132").Value = 206
$ws.Range("T132").Value = 51
$ws.Range("U132").Value = 370
$ws.Range("V132").Value = 193
$ws.Range("W132").Value = 476
$ws.Range("X132").Value = 566
$ws.Range("Y132").Value = 68
$ws.Range("Z132").Value = 496
$ws.Range("AA132").Value = 257

# --- New rows: 2022 week 27 and week 28 (monkeypox update period) ---
$ws.Range("P133").Value = 95
$ws.Range("Q133").Value = 131
$ws.Range("R133").Value = 116
$ws.Range("S133").Value = 183
$ws.Range("T133").Value = 52
$ws.Range("U133").Value = 379
$ws.Range("V133").Value = 194
$ws.Range("W133").Value = 467
$ws.Range("X133").Value = 597
$ws.Range("Y133").Value = 77
$ws.Range("Z133").Value = 443
$ws.Range("AA133").Value = 205
$ws.Range("AC133").Value = 2022
$ws.Range("AD133").Value = 27
$ws.Range("AE133").Formula = "=ROUND((P133-B133)/B133*100,2)"
$ws.Range("AF133").Formula = "=ROUND((Q133-C133)/C133*100,2)"
$ws.Range("AG133").Formula = "=ROUND((R133-D133)/D133*100,2)"
$ws.Range("AH133").Formula = "=ROUND((S133-E133)/E133*100,2)"
$ws.Range("AI133").Formula = "=ROUND((T133-F133)/F133*100,2)"
$ws.Range("AJ133").Formula = "=ROUND((U133-G133)/G133*100,2)"
$ws.Range("AK133").Formula = "=ROUND((V133-H133)/H133*100,2)"
$ws.Range("AL133").Formula = "=ROUND((W133-I133)/I133*100,2)"
$ws.Range("AM133").Formula = "=ROUND((X133-J133)/J133*100,2)"
$ws.Range("AN133").Formula = "=ROUND((Y133-K133)/K133*100,2)"
$ws.Range("AO133").Formula = "=ROUND((Z133-L133)/L133*100,2)"
$ws.Range("AP133").Formula = "=ROUND((AA133-M133)/M133*100,2)"
$ws.Range("P134").Value = 124
$ws.Range("Q134").Value = 118
$ws.Range("R134").Value = 111
$ws.Range("S134").Value = 196
$ws.Range("T134").Value = 58
$ws.Range("U134").Value = 349
$ws.Range("V134").Value = 219
$ws.Range("W134").Value = 529
$ws.Range("X134").Value = 650
$ws.Range("Y134").Value = 79
$ws.Range("Z134").Value = 440
$ws.Range("AA134").Value = 264
$ws.Range("AC134").Value = 2022
$ws.Range("AD134").Value = 28
$ws.Range("AE134").Formula = "=ROUND((P134-B134)/B134*100,2)"
$ws.Range("AF134").Formula = "=ROUND((Q134-C134)/C134*100,2)"
$ws.Range("AG134").Formula = "=ROUND((R134-D134)/D134*100,2)"
$ws.Range("AH134").Formula = "=ROUND((S134-E134)/E134*100,2)"
$ws.Range("AI134").Formula = "=ROUND((T134-F134)/F134*100,2)"
$ws.Range("AJ134").Formula = "=ROUND((U134-G134)/G134*100,2)"
$ws.Range("AK134").Formula = "=ROUND((V134-H134)/H134*100,2)"
$ws.Range("AL134").Formula = "=ROUND((W134-I134)/I134*100,2)"
$ws.Range("AM134").Formula = "=ROUND((X134-J134)/J134*100,2)"
$ws.Range("AN134").Formula = "=ROUND((Y134-K134)/K134*100,2)"
$ws.Range("AO134").Formula = "=ROUND((Z134-L134)/L134*100,2)"
$ws.Range("AP134").Formula = "=ROUND((AA134-M134)/M134*100,2)"

# --- Leave the sheet selection where the author left it after the edit ---
$ws.Range("AI133").Select()
